$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "40.941.43"
$ws.Range("E2").Value = "  -1.77%  "

# Row 3
$ws.Range("D3").Value = "2.414.94"
$ws.Range("E3").Value = "  -2.41%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.80"
$ws.Range("E5").Value = "  -1.26%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.31"
$ws.Range("E6").Value = "  -4.98%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.536"
$ws.Range("E7").Value = "  -2.86%  "

# Row 8
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  -4.30%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0832"
$ws.Range("E10").Value = "  -2.70%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.32"
$ws.Range("E11").Value = "  -5.48%  "

# Row 12
$ws.Range("E12").Value = "  -1.70%  "

# Row 13
$ws.Range("D13").Value = "2.788.02"
$ws.Range("E13").Value = "  -2.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.75"
$ws.Range("E14").Value = "  -2.17%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.52"
$ws.Range("E15").Value = "  -1.69%  "

# Row 16
$ws.Range("D16").Value = "2.416.29"
$ws.Range("E16").Value = "  -2.18%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.768"
$ws.Range("E17").Value = "  -2.63%  "

# Row 18
$ws.Range("D18").Value = "40.862.26"
$ws.Range("E18").Value = "  -1.86%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0918"
$ws.Range("E19").Value = "  -3.51%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.22"
$ws.Range("E20").Value = "  -3.94%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.54"
$ws.Range("E21").Value = "  -1.04%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.83"
$ws.Range("E22").Value = "  -4.23%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.49"
$ws.Range("E23").Value = "  -0.47%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.66"
$ws.Range("E24").Value = "  -3.17%  "

# Row 25
$ws.Range("E25").Value = "  +0.13%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.84"
$ws.Range("E26").Value = "  -4.58%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.94"
$ws.Range("E27").Value = "  -3.43%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -2.62%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.48"
$ws.Range("E29").Value = "  -3.72%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.05"
$ws.Range("E30").Value = "  -5.65%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.23"
$ws.Range("E31").Value = "  -1.03%  "

# Row 32
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.17%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.25"
$ws.Range("E33").Value = "  -4.95%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0735"
$ws.Range("E34").Value = "  -4.28%  "

# Row 35
$ws.Range("E35").Value = "  -4.85%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.86"
$ws.Range("E36").Value = "  -2.45%  "

# Row 37
$ws.Range("E37").Value = "  -1.41%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.15"
$ws.Range("E38").Value = "  -6.96%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.75"
$ws.Range("E39").Value = "  -6.74%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0997"
$ws.Range("E40").Value = "  -3.75%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.85"
$ws.Range("E41").Value = "  -3.82%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.29"
$ws.Range("E42").Value = "  -6.85%  "

# Row 43
$ws.Range("D43").Value = "1.984.88"
$ws.Range("E43").Value = "  -0.61%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.22"
$ws.Range("E44").Value = "  -4.53%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0272"
$ws.Range("E45").Value = "  -4.49%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.84"
$ws.Range("E46").Value = "  -4.74%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.36"
$ws.Range("E47").Value = "  +0.35%  "

# Row 48
$ws.Range("D48").Value = "2.646.89"
$ws.Range("E48").Value = "  -2.42%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "73.66"
$ws.Range("E49").Value = "  -0.73%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "93.81"
$ws.Range("E50").Value = "  -3.67%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.23"
$ws.Range("E51").Value = "  -2.07%  "
